$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price text would otherwise be auto-parsed as a number by Excel;
# force Text format first so the literal string (incl. trailing zeros) is preserved.
$textForceCells = @("D5", "D6", "D9", "D10", "D12", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D27", "D28", "D30", "D31", "D33", "D35", "D37", "D40", "D41", "D42", "D43", "D46", "D47", "D49")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "63.078.74"
$ws.Range("E2").Value = "  +0.61%  "

$ws.Range("D3").Value = "3.070.51"
$ws.Range("E3").Value = "  +0.32%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").Value = "538.69"
$ws.Range("E5").Value = "  -0.16%  "

$ws.Range("D6").Value = "137.23"
$ws.Range("E6").Value = "  +2.70%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.065.83"
$ws.Range("E8").Value = "  +0.30%  "

$ws.Range("D9").Value = "0.493"
$ws.Range("E9").Value = "  +0.87%  "

$ws.Range("D10").Value = "0.156"
$ws.Range("E10").Value = "  +0.83%  "

$ws.Range("E11").Value = "  +1.77%  "

$ws.Range("D12").Value = "0.454"
$ws.Range("E12").Value = "  -1.82%  "

$ws.Range("E13").Value = "  +0.43%  "

$ws.Range("D14").Value = "34.44"
$ws.Range("E14").Value = "  -0.88%  "

$ws.Range("D15").Value = "3.562.46"
$ws.Range("E15").Value = "  +0.15%  "

$ws.Range("D16").Value = "63.013.12"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("E17").Value = "  +1.63%  "

$ws.Range("D18").Value = "3.067.14"
$ws.Range("E18").Value = "  +0.19%  "

$ws.Range("D19").Value = "6.63"
$ws.Range("E19").Value = "  -0.82%  "

$ws.Range("D20").Value = "469.80"
$ws.Range("E20").Value = "  -1.82%  "

$ws.Range("D21").Value = "13.50"
$ws.Range("E21").Value = "  +0.45%  "

$ws.Range("D22").Value = "0.695"
$ws.Range("E22").Value = "  -1.93%  "

$ws.Range("D23").Value = "7.02"
$ws.Range("E23").Value = "  -2.46%  "

$ws.Range("D24").Value = "78.44"
$ws.Range("E24").Value = "  -0.09%  "

$ws.Range("D25").Value = "12.11"
$ws.Range("E25").Value = "  +0.66%  "

$ws.Range("E26").Value = "  -0.08%  "

$ws.Range("D27").Value = "2.69"
$ws.Range("E27").Value = "  -0.82%  "

$ws.Range("D28").Value = "7.86"
$ws.Range("E28").Value = "  -4.88%  "

$ws.Range("E29").Value = "  -0.12%  "

$ws.Range("D30").Value = "26.12"
$ws.Range("E30").Value = "  -0.20%  "

$ws.Range("D31").Value = "1.16"
$ws.Range("E31").Value = "  +4.99%  "

$ws.Range("E32").Value = "  -2.67%  "

$ws.Range("D33").Value = "58.97"
$ws.Range("E33").Value = "  +1.53%  "

$ws.Range("E34").Value = "  -5.22%  "

$ws.Range("D35").Value = "5.45"
$ws.Range("E35").Value = "  +6.79%  "

$ws.Range("E36").Value = "  -0.22%  "

$ws.Range("D37").Value = "480.31"
$ws.Range("E37").Value = "  -2.12%  "

$ws.Range("D38").Value = "3.254.65"
$ws.Range("E38").Value = "  +3.55%  "

$ws.Range("E39").Value = "  +1.04%  "

$ws.Range("D40").Value = "0.0792"
$ws.Range("E40").Value = "  -0.71%  "

$ws.Range("D41").Value = "0.117"
$ws.Range("E41").Value = "  -0.11%  "

$ws.Range("D42").Value = "8.13"
$ws.Range("E42").Value = "  +0.71%  "

$ws.Range("D43").Value = "2.57"
$ws.Range("E43").Value = "  +1.27%  "

$ws.Range("E44").Value = "  -0.89%  "

$ws.Range("E45").Value = "  +0.12%  "

$ws.Range("D46").Value = "123.52"
$ws.Range("E46").Value = "  +4.51%  "

$ws.Range("D47").Value = "25.22"
$ws.Range("E47").Value = "  +2.19%  "

$ws.Range("E48").Value = "  -1.58%  "

$ws.Range("D49").Value = "0.109"
$ws.Range("E49").Value = "  +0.96%  "

$ws.Range("E50").Value = "  +2.37%  "

$ws.Range("E51").Value = "  -0.20%  "
